$p = $ppt.ActivePresentation

# The table on slide 5 (the "Type of document / Definition / Why it is
# important" grid) had its table style switched away from the deck's
# custom style to a different (built-in) table style. Locate that table
# - it's the 2nd shape on slide 5 - and re-apply the new style's GUID.
$s = $p.Slides.Item(5)
$sh = $s.Shapes.Item(2)

if ($sh.HasTable) {
    $tbl = $sh.Table
    $tbl.ApplyStyle("{EB5382F3-0EC5-4EFB-8EEC-55ADECB99FF4}")
}
